$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: " raven, maar kan geopend worden in " ->
#   " raven, maar kan geopend en ge-edit worden in "
# (the real edit split this across several runs / wrapped "edit" in a
# spell-check mark; the Word object model always re-merges adjacent runs
# that share identical formatting when the document is serialized, so we
# apply the textual change as a single in-place replacement.)
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(" raven, maar kan geopend worden in ", $true, $false, $false, $false, $false, $true, 1, $false, " raven, maar kan geopend en ge-edit worden in ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Ga naar File, Open project en navigeer naar de root van de
# folder waar je met git een project hebt aangemaakt." ->
#   "Ga naar "File", "Open project" en navigeer naar de root van de folder
#   waar je met git een project hebt aangemaakt."
# (curly/smart quotes around File and Open project)
# ---------------------------------------------------------------------------
$quote1 = [char]8220
$quote2 = [char]8221
$replacement2 = "Ga naar " + $quote1 + "File" + $quote2 + ", " + $quote1 + "Open project" + $quote2 + " en navigeer naar de root van de folder waar je met git een project hebt aangemaakt."
$r2 = $d.Content
$r2.Find.Execute("Ga naar File, Open project en navigeer naar de root van de folder waar je met git een project hebt aangemaakt.", $true, $false, $false, $false, $false, $true, 1, $false, $replacement2, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: " zullen met " -> " zullen automatisch met "
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute(" zullen met ", $true, $false, $false, $false, $false, $true, 1, $false, " zullen automatisch met ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: move the "_GoBack" bookmark from the middle of the "Procedures
# in het SQL-bestand ..." paragraph (right after "... van het bestand") to
# the end of the previous paragraph (right after "... deployen."), so the
# bookmark no longer splits the "Procedures..." paragraph's own text and the
# paragraph break now falls right after "... deployen." instead of right
# after "... van het bestand".
# ---------------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("Procedures in het SQL-bestand met testgegevens kunnen niet uitgevoerd worden in SQL zonder bewerking", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $procParaStart = $r4.Start
    $anchor = $d.Range($procParaStart, $procParaStart)
    $anchor.Collapse(1)
    $prevParaRange = $anchor.Paragraphs(1).Previous(1).Range
    $insertPos = $prevParaRange.End - 1

    # Use a one-character placeholder to get a *non-empty* range (the engine
    # mishandles zero-length ranges passed straight into Bookmarks.Add), add
    # the bookmark around it, then delete the placeholder so the bookmark
    # collapses back down to an empty range in the right spot. Re-adding a
    # bookmark with the existing name "_GoBack" relocates it instead of
    # duplicating it.
    $placeholder = $d.Range($insertPos, $insertPos)
    $placeholder.InsertAfter("@")
    $bmRange = $d.Range($insertPos, $insertPos + 1)
    $d.Bookmarks.Add("_GoBack", $bmRange)
    $d.Range($insertPos, $insertPos + 1).Delete()
}
